# Remove the "Fifthswhy" hero row from the heroes table.
# Deleting the entire row shifts subsequent rows up, shrinks the
# table/dimension ranges automatically, and updates the shared-strings
# table via the engine's own save logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Rows.Item(13).Delete()

# Leave the selection on the row that used to be the deleted one,
# matching the state Excel leaves behind after a row delete.
$ws.Rows.Item(13).Select()
